# Weekly update: a new Perejil price record was collected for
# Terminal Hortofrutícola Agro Chillán and inserted as the new row 47,
# pushing all subsequent records (old rows 47-98) down by one row
# (new rows 48-99). Dimension grows from A1:R98 to A1:R99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 47; everything below shifts down one row.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new record.
$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C47").Value = "Ñuble"
$ws.Range("D47").Value = 45117
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = 100112044
$ws.Range("G47").Value = "Perejil"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 1500
$ws.Range("L47").Value = 1500
$ws.Range("M47").Value = 1500
$ws.Range("N47").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O47").Value = "Región de Ñuble"
$ws.Range("P47").Value = 1500
$ws.Range("Q47").Value = 1
$ws.Range("R47").Value = "Hortaliza"
